$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.343.22'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '1.689.00'
$ws.Range("E3").Value = '  +1.32%  '

$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").Value = '218.45'
$ws.Range("E5").Value = '  +0.51%  '

$ws.Range("D6").Value = '0.5260'
$ws.Range("E6").Value = '  +3.80%  '

$ws.Range("D7").Value = '1.008'
$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("D8").Value = '0.2706'
$ws.Range("E8").Value = '  +2.14%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.06423'
$ws.Range("E9").Value = '  +1.26%  '

$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").Value = '22.05'
$ws.Range("E10").Value = '  +2.92%  '

$ws.Range("D11").Value = '0.07514'
$ws.Range("E11").Value = '  +2.06%  '

$ws.Range("D12").Value = '1.718.71'
$ws.Range("E12").Value = '  +2.97%  '

$ws.Range("D13").Value = '4.574'
$ws.Range("E13").Value = '  +0.85%  '

$ws.Range("D14").Value = '0.5830'
$ws.Range("E14").Value = '  +0.50%  '

$ws.Range("D15").Value = '0.000008515'
$ws.Range("E15").Value = '  +0.25%  '

$ws.Range("D16").Value = '64.54'
$ws.Range("E16").Value = '  -0.04%  '

$ws.Range("D17").Value = '26.385.55'
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").Value = '4.937'
$ws.Range("E18").Value = '  +0.46%  '

$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").Value = '10.90'
$ws.Range("E20").Value = '  +0.85%  '

$ws.Range("D21").Value = '189.36'
$ws.Range("E21").Value = '  +0.55%  '

$ws.Range("D22").Value = '6.217'
$ws.Range("E22").Value = '  +0.74%  '

$ws.Range("D23").Value = '1.008'
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").Value = '144.62'
$ws.Range("E24").Value = '  +1.10%  '

$ws.Range("D25").Value = '7.736'
$ws.Range("E25").Value = '  +0.99%  '

$ws.Range("D26").Value = '0.1233'
$ws.Range("E26").Value = '  +5.42%  '

$ws.Range("E27").Value = '  +1.74%  '

$ws.Range("D28").Value = '0.06635'
$ws.Range("E28").Value = '  +14.03%  '

$ws.Range("D29").Value = '1.353'
$ws.Range("E29").Value = '  +6.90%  '


$ws.Range("D31").Value = '3.584'
$ws.Range("E31").Value = '  +2.28%  '

$ws.Range("D32").Value = '3.573'
$ws.Range("E32").Value = '  +1.44%  '

$ws.Range("D33").Value = '1.664'
$ws.Range("E33").Value = '  +1.90%  '

$ws.Range("D34").Value = '1.029'
$ws.Range("E34").Value = '  +1.90%  '

$ws.Range("D35").Value = '0.6236'
$ws.Range("E35").Value = '  +4.25%  '

$ws.Range("D36").Value = '2.396'
$ws.Range("E36").Value = '  +1.57%  '

$ws.Range("D37").Value = '2.707'
$ws.Range("E37").Value = '  +2.33%  '

$ws.Range("D38").Value = '6.377'
$ws.Range("E38").Value = '  +6.21%  '

$ws.Range("D39").Value = '1.114.89'
$ws.Range("E39").Value = '  +3.95%  '

$ws.Range("D40").Value = '0.01619'
$ws.Range("E40").Value = '  +0.69%  '

$ws.Range("D41").Value = '0.8865'
$ws.Range("E41").Value = '  +2.69%  '

$ws.Range("E42").Value = '  +0.85%  '

$ws.Range("D43").Value = '101.10'
$ws.Range("E43").Value = '  +1.52%  '

$ws.Range("D44").Value = '1.838.24'
$ws.Range("E44").Value = '  +1.26%  '

$ws.Range("D45").Value = '0.00000000115'
$ws.Range("E45").Value = '  +4.01%  '

$ws.Range("D46").Value = '56.94'
$ws.Range("E46").Value = '  +2.39%  '

$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '1.008'
$ws.Range("E47").Value = '  +0.42%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '8.149'
$ws.Range("E48").Value = '  +0.97%  '

$ws.Range("D49").Value = '0.05270'
$ws.Range("E49").Value = '  +1.80%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '0.4306'
$ws.Range("E50").Value = '  +0.29%  '

$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '6.074'
$ws.Range("E51").Value = '  +3.99%  '
